# 2.a.1.xlsx update: add a new "2021" data column (column R) to the table.
#
# - Column R gets a header value of 2021 in row 3 (same look as the
#   existing year headers in O3:Q3).
# - Column R gets the new data value in row 4 (same look as the existing
#   data cells in D4:Q4), except Excel created a brand-new font/style
#   record for it (a duplicate of the one used by D4:Q4) - we reproduce
#   that by nudging the cell's font so a new font/style entry is created
#   instead of reusing the existing one.
# - The two preceding values (P4, Q4) were recalculated/updated.
# - The workbook was left with cell O10 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell R3 (2021), formatted like O3:Q3 ---
$ws.Range("Q3").Copy() | Out-Null
$ws.Range("R3").PasteSpecial(-4122) | Out-Null
$ws.Range("R3").Value = 2021

# --- Update the two existing percentages that changed ---
$ws.Range("P4").Value = 0.09130340807234763
$ws.Range("Q4").Value = 0.074862480994528399

# --- New data cell R4, formatted like D4:Q4 ---
$ws.Range("Q4").Copy() | Out-Null
$ws.Range("R4").PasteSpecial(-4122) | Out-Null
$ws.Range("R4").Value = 0.064467421337540437

# Force a distinct (new) font/style record for R4, matching how the
# published workbook ends up with an extra font + cellXf pair applied
# only to this cell.
$ws.Range("R4").Font.ThemeFont = 1

# Clear the clipboard marching ants / copy mode
$excel.CutCopyMode = 0

# --- Leave the selection where the author left it ---
$ws.Range("O10").Select() | Out-Null
